# Add a new "2022-Q3" sheet of quarterly fund-holding data for this stock,
# inserted right after the "总计" (totals) sheet, and update the "总计"
# sheet to include this new quarter in its summary table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet immediately before "2022-Q2"
#    (i.e. right after "总计"), matching the tab order produced when a
#    new quarter is inserted ahead of the existing history. Duplicating
#    the existing "2022-Q2" sheet keeps the same layout/formatting
#    (column styles, borders, page margins) used by all the other
#    quarterly sheets.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)
$template.Copy($template)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template only has 5 data rows (rows 2-6); the new quarter needs 7
# (rows 2-8), so extend the table by copying the formatting of the last
# data row down two more rows.
$q3.Range("A6:H6").Copy()
$q3.Range("A7:H7").PasteSpecial(-4122)
$q3.Range("A6:H6").Copy()
$q3.Range("A8:H8").PasteSpecial(-4122)

# Header row (kept the same text as the template, rewritten explicitly
# for clarity/safety)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Force columns B:G on the data rows to be plain text so that codes like
# "512980" and decimal strings like "44.76" are preserved verbatim
# instead of being auto-converted to numbers.
$q3.Range("B2:G8").NumberFormat = "@"

$q3Data = @(
    @("512980", "广发中证传媒ETF",               "44.76", "99.29", "3.15", "1.4099", 7),
    @("160629", "鹏华中证传媒指数（LOF）A",        "6.41",  "94.58", "2.97", "0.1904", 7),
    @("159805", "鹏华中证传媒ETF",                 "1.71",  "98.37", "3.11", "0.0532", 7),
    @("164818", "工银瑞信中证传媒指数（LOF）A",     "1.65",  "93.46", "2.94", "0.0485", 7),
    @("010677", "工银瑞信中证传媒指数（LOF）C",     "0.21",  "93.46", "2.94", "0.0062", 7),
    @("015675", "鹏华中证传媒指数（LOF）C",        "0.17",  "94.58", "2.97", "0.0050", 7),
    @("516190", "华夏中证文娱传媒ETF",             "0.13",  "96.01", "2.56", "0.0033", 10)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: a new top row for 2022-Q3 is
#    inserted and the rest of the quarterly history shifts down.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item(1)

$zjData = @(
    @("2022-Q3", 7, 1.72),
    @("2022-Q2", 5, 1.72),
    @("2022-Q1", 1, 0.24),
    @("2021-Q4", 1, 0.27),
    @("2021-Q3", 5, 2.2),
    @("2021-Q2", 5, 1.42),
    @("2021-Q1", 8, 1.74),
    @("2020-Q4", 6, 1.36)
)

$r = 2
foreach ($row in $zjData) {
    $zj.Cells.Item($r, 1).Value = $r - 2
    $zj.Cells.Item($r, 2).Value = $row[0]
    $zj.Cells.Item($r, 3).Value = $row[1]
    $zj.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# The new last row (A9) needs the same "index column" formatting as the
# rest of column A; copy it down from the row above.
$zj.Range("A8").Copy()
$zj.Range("A9").PasteSpecial(-4122)
